# Auto-generated edit script applying the Durandal_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 980.2
$ws.Range("I4").Value = 675.25
$ws.Range("J4").Value = 2200
$ws.Range("K4").Value = 675.25
$ws.Range("L4").Value = 2200
$ws.Range("M4").Value = -561.25
$ws.Range("N4").Value = -2428
$ws.Range("H64").Value = 3747
$ws.Range("I64").Value = 2990
$ws.Range("J64").Value = 4071.4285
$ws.Range("K64").Value = 2990
$ws.Range("L64").Value = 4071.4285
$ws.Range("M64").Value = -2742
$ws.Range("N64").Value = -4567.4285
$ws.Range("H67").Value = 3747
$ws.Range("I67").Value = 2990
$ws.Range("J67").Value = 4071.4285
$ws.Range("K67").Value = 2990
$ws.Range("L67").Value = 4071.4285
$ws.Range("M67").Value = -2132
$ws.Range("N67").Value = -5787.4285
$ws.Range("H98").Value = 55563972
$ws.Range("I98").Value = 9984.25
$ws.Range("K98").Value = 9984.25
$ws.Range("M98").Value = -8486.25
$ws.Range("H100").Value = 10754372
$ws.Range("I100").Value = 12346824
$ws.Range("J100").Value = 5321.5
$ws.Range("K100").Value = 12346824
$ws.Range("L100").Value = 5321.5
$ws.Range("M100").Value = -12346283
$ws.Range("N100").Value = -6403.5
$ws.Range("H113").Value = 2072.0605
$ws.Range("I113").Value = 2112.7273
$ws.Range("K113").Value = 2112.7273
$ws.Range("M113").Value = 1141.2727
$ws.Range("H122").Value = 55563972
$ws.Range("I122").Value = 9984.25
$ws.Range("K122").Value = 29952.75
$ws.Range("M122").Value = -27502.75
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
$ws.Range("H137").Value = 1162.5714
$ws.Range("I137").Value = 988
$ws.Range("J137").Value = 1511.7142
$ws.Range("K137").Value = 2964
$ws.Range("L137").Value = 4535.142599999999
$ws.Range("M137").Value = -414
$ws.Range("N137").Value = -9635.142599999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2812
$ws.Range("I45").Value = 2676.9473
$ws.Range("J45").Value = 3045.2727
$ws.Range("K45").Value = 2676.9473
$ws.Range("L45").Value = 3045.2727
$ws.Range("M45").Value = -2299.9473
$ws.Range("N45").Value = -3799.2727
$ws.Range("H61").Value = 2175306.8
$ws.Range("I61").Value = 2632660.8
$ws.Range("J61").Value = 2875
$ws.Range("K61").Value = 2632660.8
$ws.Range("L61").Value = 2875
$ws.Range("M61").Value = -2632448.8
$ws.Range("N61").Value = -3299
$ws.Range("H74").Value = 635.88464
$ws.Range("I74").Value = 633.5625
$ws.Range("J74").Value = 639.6
$ws.Range("K74").Value = 633.5625
$ws.Range("L74").Value = 639.6
$ws.Range("M74").Value = 240.4375
$ws.Range("N74").Value = -2387.6
$ws.Range("H77").Value = 635.88464
$ws.Range("I77").Value = 633.5625
$ws.Range("J77").Value = 639.6
$ws.Range("K77").Value = 3167.8125
$ws.Range("L77").Value = 3198
$ws.Range("M77").Value = 1200.1875
$ws.Range("N77").Value = -11934
$ws.Range("H136").Value = 2175306.8
$ws.Range("I136").Value = 2632660.8
$ws.Range("J136").Value = 2875
$ws.Range("K136").Value = 7897982.399999999
$ws.Range("L136").Value = 8625
$ws.Range("M136").Value = -7895432.399999999
$ws.Range("N136").Value = -13725

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2292.7454
$ws.Range("J31").Value = 4213.3335
$ws.Range("L31").Value = 4213.3335
$ws.Range("N31").Value = -4803.3335
$ws.Range("H34").Value = 2292.7454
$ws.Range("J34").Value = 4213.3335
$ws.Range("L34").Value = 4213.3335
$ws.Range("N34").Value = -4617.3335
$ws.Range("H64").Value = 10000
$ws.Range("J64").Value = 10000
$ws.Range("L64").Value = 10000
$ws.Range("N64").Value = -10496
$ws.Range("H67").Value = 10000
$ws.Range("J67").Value = 10000
$ws.Range("L67").Value = 10000
$ws.Range("N67").Value = -11716
$ws.Range("H86").Value = 131712.86
$ws.Range("J86").Value = 136998.33
$ws.Range("L86").Value = 136998.33
$ws.Range("N86").Value = -139244.33
$ws.Range("H89").Value = 131712.86
$ws.Range("J89").Value = 136998.33
$ws.Range("L89").Value = 684991.6499999999
$ws.Range("N89").Value = -696223.6499999999
$ws.Range("H122").Value = 1800
$ws.Range("I122").Value = 1766.6666
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 5299.9998
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -2849.9998
$ws.Range("N122").Value = -10900
$ws.Range("H132").Value = 74613.71000000001
$ws.Range("I132").Value = 2445.3333
$ws.Range("J132").Value = 128740
$ws.Range("K132").Value = 7335.999899999999
$ws.Range("L132").Value = 386220
$ws.Range("M132").Value = -4805.999899999999
$ws.Range("N132").Value = -391280

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 334233.34
$ws.Range("I117").Value = 1350
$ws.Range("K117").Value = 4050
$ws.Range("M117").Value = -608
$ws.Range("H122").Value = 600.8889
$ws.Range("I122").Value = 334.66666
$ws.Range("J122").Value = 734
$ws.Range("K122").Value = 3011.99994
$ws.Range("L122").Value = 6606
$ws.Range("M122").Value = -561.9999399999997
$ws.Range("N122").Value = -11506
$ws.Range("H125").Value = 5600
$ws.Range("I125").Value = 1000
$ws.Range("J125").Value = 6750
$ws.Range("K125").Value = 3000
$ws.Range("L125").Value = 20250
$ws.Range("M125").Value = 1920
$ws.Range("N125").Value = -30090
$ws.Range("H131").Value = 6250949
$ws.Range("J131").Value = 7043200
$ws.Range("L131").Value = 21129600
$ws.Range("N131").Value = -21139680
$ws.Range("H137").Value = 2107.111
$ws.Range("J137").Value = 3493.2222
$ws.Range("L137").Value = 10479.6666
$ws.Range("N137").Value = -20679.6666
$ws.Range("H139").Value = 5100.6665
$ws.Range("I139").Value = 5537.273
$ws.Range("J139").Value = 3900
$ws.Range("K139").Value = 16611.819
$ws.Range("L139").Value = 11700
$ws.Range("M139").Value = -11471.819
$ws.Range("N139").Value = -21980

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4632.857
$ws.Range("I80").Value = 4725
$ws.Range("J80").Value = 4596
$ws.Range("K80").Value = 4725
$ws.Range("L80").Value = 4596
$ws.Range("M80").Value = -3727
$ws.Range("N80").Value = -6592
$ws.Range("H83").Value = 4632.857
$ws.Range("I83").Value = 4725
$ws.Range("J83").Value = 4596
$ws.Range("K83").Value = 23625
$ws.Range("L83").Value = 22980
$ws.Range("M83").Value = -18633
$ws.Range("N83").Value = -32964
$ws.Range("H102").Value = 1632.2354
$ws.Range("I102").Value = 1648.9333
$ws.Range("J102").Value = 1507
$ws.Range("K102").Value = 1648.9333
$ws.Range("L102").Value = 1507
$ws.Range("M102").Value = -26.93329999999992
$ws.Range("N102").Value = -4751
$ws.Range("H107").Value = 396.85715
$ws.Range("I107").Value = 251
$ws.Range("J107").Value = 542.7143
$ws.Range("K107").Value = 251
$ws.Range("L107").Value = 542.7143
$ws.Range("M107").Value = 1669
$ws.Range("N107").Value = -4382.7143
$ws.Range("H122").Value = 1776.2222
$ws.Range("I122").Value = 1722
$ws.Range("J122").Value = 1884.6666
$ws.Range("K122").Value = 5166
$ws.Range("L122").Value = 5653.9998
$ws.Range("M122").Value = -2716
$ws.Range("N122").Value = -10553.9998
$ws.Range("H126").Value = 20834300
$ws.Range("I126").Value = 1064.2
$ws.Range("K126").Value = 3192.6
$ws.Range("M126").Value = -722.6000000000004

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2261.5757
$ws.Range("I7").Value = 1728.5927
$ws.Range("K7").Value = 1728.5927
$ws.Range("M7").Value = -1616.5927
$ws.Range("H40").Value = 1853.1538
$ws.Range("I40").Value = 1924.5
$ws.Range("J40").Value = 1739
$ws.Range("K40").Value = 1924.5
$ws.Range("L40").Value = 1739
$ws.Range("M40").Value = -1788.5
$ws.Range("N40").Value = -2011
$ws.Range("H46").Value = 9417.166999999999
$ws.Range("I46").Value = 920.8
$ws.Range("J46").Value = 15486
$ws.Range("K46").Value = 920.8
$ws.Range("L46").Value = 15486
$ws.Range("M46").Value = -732.8
$ws.Range("N46").Value = -15862
$ws.Range("H61").Value = 3199.75
$ws.Range("I61").Value = 2399.5
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2399.5
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2197.5
$ws.Range("N61").Value = -4404
$ws.Range("H82").Value = 2871.6667
$ws.Range("I82").Value = 3576.6667
$ws.Range("J82").Value = 2166.6667
$ws.Range("K82").Value = 3576.6667
$ws.Range("L82").Value = 2166.6667
$ws.Range("M82").Value = -3215.6667
$ws.Range("N82").Value = -2888.6667
$ws.Range("H85").Value = 2871.6667
$ws.Range("I85").Value = 3576.6667
$ws.Range("J85").Value = 2166.6667
$ws.Range("K85").Value = 3576.6667
$ws.Range("L85").Value = 2166.6667
$ws.Range("M85").Value = -2328.6667
$ws.Range("N85").Value = -4662.6667
$ws.Range("H113").Value = 3199.75
$ws.Range("I113").Value = 2399.5
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2399.5
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -229.5
$ws.Range("N113").Value = -8340
$ws.Range("H122").Value = 1947.1389
$ws.Range("J122").Value = 2061.4614
$ws.Range("L122").Value = 6184.3842
$ws.Range("N122").Value = -11084.3842
$ws.Range("H126").Value = 2261.5757
$ws.Range("I126").Value = 1728.5927
$ws.Range("K126").Value = 5185.7781
$ws.Range("M126").Value = -2715.7781
$ws.Range("H132").Value = 6360.0835
$ws.Range("I132").Value = 7299.107
$ws.Range("J132").Value = 3073.5
$ws.Range("K132").Value = 21897.321
$ws.Range("L132").Value = 9220.5
$ws.Range("M132").Value = -19367.321
$ws.Range("N132").Value = -14280.5
$ws.Range("H136").Value = 9687.611000000001
$ws.Range("I136").Value = 8676.571
$ws.Range("K136").Value = 26029.713
$ws.Range("M136").Value = -23479.713

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2577.5
$ws.Range("I122").Value = 2710.7144
$ws.Range("J122").Value = 2266.6667
$ws.Range("K122").Value = 8132.1432
$ws.Range("L122").Value = 6800.000100000001
$ws.Range("M122").Value = -5682.1432
$ws.Range("N122").Value = -11700.0001
$ws.Range("H132").Value = 38864240
$ws.Range("I132").Value = 59530116
$ws.Range("J132").Value = 2698965.2
$ws.Range("K132").Value = 178590348
$ws.Range("L132").Value = 8096895.600000001
$ws.Range("M132").Value = -178587818
$ws.Range("N132").Value = -8101955.600000001
